$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Transactions
# Delete the four "Accrual" transaction rows (rows 2, 3, 5, 7), leaving
# the disburse/repayment rows which shift up to rows 2-4. Then patch the
# remaining rows with their updated figures.
# ---------------------------------------------------------------------
$wsTx = $wb.Worksheets.Item("Transactions")
$wsTx.Activate()

$wsTx.Rows.Item(7).Delete()
$wsTx.Rows.Item(5).Delete()
$wsTx.Rows.Item(3).Delete()
$wsTx.Rows.Item(2).Delete()

$wsTx.Range("A2").Value = 71
$wsTx.Range("J2").Value = 9163.24
$wsTx.Range("J2").NumberFormat = "#,##0.00"

$wsTx.Range("A3").Value = 69
$wsTx.Range("F3").Value = 836.76
$wsTx.Range("H3").Value = 50
$wsTx.Range("J3").Value = 4163.24
$wsTx.Range("J3").NumberFormat = "#,##0.00"

$wsTx.Range("A4").Value = 66

$wsTx.Range("A2:XFD9").Select()

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()

$wsSummary.Range("B2").Value = 836.76
$wsSummary.Range("E2").Value = 9163.24
$wsSummary.Range("F2").Value = 877.35
$wsSummary.Range("A4").Value = 50
$wsSummary.Range("B4").Value = 50

$wsSummary.Range("B4").Select()

# ---------------------------------------------------------------------
# Sheet: Repayment schedule
# ---------------------------------------------------------------------
$wsSched = $wb.Worksheets.Item("Repayment schedule")
$wsSched.Activate()

$wsSched.Range("I3").Value = 50
$wsSched.Range("K3").Value = 937.72
$wsSched.Range("P3").Value = 27.95

$wsSched.Range("H21").Select()

# ---------------------------------------------------------------------
# Leave the workbook with the Transactions tab active, matching the
# original file's last-active-sheet state.
# ---------------------------------------------------------------------
$wsTx.Activate()
